$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Accountant- Family Office"
$ws.Range("B2").Value = "The Quest Organization"
$ws.Range("C2").Value = "New York, NY"
$ws.Range("D2").Value = "Full-Time"
$ws.Range("E2").Value = "Strong knowledge of ground-up stick frame (podium, wrap, garden) multifamily projects * 6+ years construction experience as a Construction Manager or Superintendent * Experience with managing ..."
$ws.Range("F2").Value = "https://www.ziprecruiter.com/k/l/AAK5bH-ujX0aetfDLmO_NKbdWMsbwFlxTXv7Mm_HbSLUWqD3vwI5avQ5BoK4tKTSg9GhKaBQWpj2P2S0gM8HDLk0_vewypxD7LRhzm-Ztn7DQuokXD69IR8Bj0hGn3bx-u6APq5XqpHnms9WT-y6J2Owm4kDIS2rrPid4ukhvuP7HG-2MTgpOA"

$ws.Range("A3").Value = "Property Accountant"
$ws.Range("B3").Value = "Creative Financial Staffing (CFS)"
$ws.Range("C3").Value = "Manhattan, NY"
$ws.Range("D3").Value = "Full-Time"
$ws.Range("E3").Value = "We are seeking a Commercial Construction Superintendent to join our team! You will oversee project planning, scheduling, subcontractors work. Responsibilities: * Oversee all aspects of construction ..."
$ws.Range("F3").Value = "https://www.ziprecruiter.com/k/l/AAIOkH8sgkTzkiHdKM3Ybh8sFPf7b4ifP9yRV4juXwH_Kcq0gjQmI-03WzV6d3z_rIlrn5Mkvq2NmyeSgOw1NAdr9vAcss7EmPLAsZMEO8mioNuj7_A2JjpOMy-7AX_1SK_4CdnwxfWP8M8H5ejL6WNa9BhTUxEclwMClhy52z8DnJUHMGcvfQ"

$ws.Range("A4").Value = "Bookkeeper Accountant"
$ws.Range("B4").Value = "Abetta Boiler Welding INC"
$ws.Range("C4").Value = "Manhattan, NY"
$ws.Range("D4").Value = "Full-Time"
$ws.Range("E4").Value = "Education * 4 yr degree in engineering, construction science, construction management or similar Position Summary We are looking for an Assistant Project Manager to perform the technical and support ..."
$ws.Range("F4").Value = "https://www.ziprecruiter.com/k/l/AAJZo6CtScvPi0RIcjnS8MJ8lEkNY12B729LLyiNf0eyPN4jWwvZ_DOrMpkDSKtfonwTDskZIt5uDp_TrPlRHyOJ2dywrQhFDOx4fJ-ZBSWW0OXZm3UTpd8qU_VDOBtbSJWIo6xJUkOdwQW35raEcsLU2KsGuEGr3A53ysMmK6F5gv8zgdBiQA"

$ws.Range("A5").Value = "Staff Accountant"
$ws.Range("B5").Value = "Indcorp Fiscal Services Inc."
$ws.Range("C5").Value = "Island Park, NY"
$ws.Range("D5").Value = "Full-Time"
$ws.Range("E5").Value = "Talent Corps is currently hiring for General Labors for the Fort Worth, TX area for a Temp to Perm position. Candidates are needing to be physically fit as this job is physically demanding and may ..."
$ws.Range("F5").Value = "https://www.ziprecruiter.com/k/l/AAKwWausp_j4tC_ejXpMvxg4uKKX1bpMnO9p1QziVwsNSl924QBlX-FjR4pinqzQ3_Y9SWmSfaT9Gt5Tuwf_151HGLYivDx9WSSWYDNUA4z3_cXl4Kb799R3G0wJQ80frPQtA9dNNTv9BFt2ff1rJp5qXIqxlbjVguePuwFxcUDBZ0u91ER-jg"

$ws.Range("A6").Value = "Junior Accountant"
$ws.Range("B6").Value = "Greater New York Insurance Companies"
$ws.Range("C6").Value = "Manhattan, NY"
$ws.Range("D6").Value = "Full-Time"
$ws.Range("E6").Value = "Have strong knowledge and experience of all phases of construction scope and sequence * Perform take-offs and estimation on plans * Enter and understand budgets * Possess experience with bid tabs (A ..."
$ws.Range("F6").Value = "https://www.ziprecruiter.com/k/l/AAJeYo2Q-gj-LOVZ8S1GbwcPnUKYqgYc3Mf6wrhgcE75tkKMScn5vGgYg-SK3VYSJpyiVxe_JYVoX2YJ6bWKhtzf57CaJSD0yQA4iSK4L6MuoqRd-o_JLMfhxtDrTFSQziuDuosUopdQGGkQJEB6hbMNbuAsJr5K97jg3F2AD99G09wOZmutvw"

$ws.Range("A7").Value = "Accountant"
$ws.Range("B7").Value = "E-Z TAX GROUP NORTH AMERICA INC"
$ws.Range("C7").Value = "Flushing, NY"
$ws.Range("D7").Value = "Full-Time"
$ws.Range("E7").Value = "Restaurant Depot is seeking a Construction & Maintenance Manager. Position will be based out of the Dallas/Richardson area offices. This person will be working under the direction of the Director of ..."
$ws.Range("F7").Value = "https://www.ziprecruiter.com/k/l/AAJURz5oMkNsaGPH5De8eLRxkY2wMDJm3_9PF3YVDWVTWaeJthqZIypT2fRqjVb93pfI97KpFo8DXRtr8FIQggnAAovU1TRaR06dsrT0A9GlbD2mgToR-4iLhL6WGW831_lduvRnS0xxkybTl2oCaTBu97c_fxX0IdYEztFU_4CDTJxCtVWdrA"

$ws.Range("A8").Value = "Staff Accountant"
$ws.Range("B8").Value = "Real Estate Management Company"
$ws.Range("C8").Value = "Bellmore, NY"
$ws.Range("D8").Value = "Full-Time"
$ws.Range("E8").Value = "Texas Home Builder seeking Construction Managers This Jobot Job is hosted by: Henry Chan Are you a fit? Easy Apply now by clicking the `"Apply Now`" button and sending us your resume. Salary: `$65,000 ..."
$ws.Range("F8").Value = "https://www.ziprecruiter.com/k/l/AAJDUWlejnI-jKXR-0-WE9cjzypZd8vj9uL2_BPJXffUfFKui8UamwypF-H4khHniLAYJYszppuvA1V-LHceD2FIWwEUh4HVEnnIFh60QW2ZZIZRngR2CJpD7LfOD1puWCADFLpBtyF-UA6RjEtv4BGJfXf82V4NrEgY5C4plkDmcjUc11mOww"

$ws.Range("A9").Value = "Senior Accountant, Ecommerce Clients"
$ws.Range("B9").Value = "Accountingfly"
$ws.Range("C9").Value = "New York, New York"
$ws.Range("D9").Value = "Full-Time"
$ws.Range("E9").Value = "The Construction Manager role is a player/coach position responsible for both managing overall construction execution of individual projects and overseeing 3rd party construction management firm in ..."
$ws.Range("F9").Value = "https://www.ziprecruiter.com/ek/l/AALZh9aBzeTQRWfyJzOuBIjByuq1vkZ0xdAKR7IxQo8FLKiN4SnXbaQ3liRLljm902fKp0Xi0T2u9Ip9bsAvp26pM4LZSVno2QsLZuLe51IQWKAZDGCK9M87_sPO6KLuI74hWCw6xrp7R2YEfvfLThbJIq0Vjm_U-C4THpRNon128EszdLvt5g"

$ws.Range("A10").Value = "Property Accountant - Real Estate"
$ws.Range("B10").Value = "Loughlin Personnel"
$ws.Range("C10").Value = "White Plains, NY"
$ws.Range("D10").Value = "Full-Time"
$ws.Range("E10").Value = "Primary Function The General Superintendent oversees the construction of multiple projects. He/she will coordinate with other superintendents to ensure other jobsites are being manned properly and ..."
$ws.Range("F10").Value = "https://www.ziprecruiter.com/k/l/AAKeCO2QEIr_QeNR19TPBBuAKVSyUimli5FAdxtLIMZD6-C5Ot3Yxo009f-UPMZNcIGR6ax2k5LRu_eaoksJNDv77UiA6D53ZADZqIA3cy7JRUc3KdRZaCk1pyi_zylOZdQ0bvirrgFLVtQCXmIu3kXeVNTYBTKXd-XagsozdZ49RF9lpI185A"

$ws.Range("A11").Value = "Staff Accountant"
$ws.Range("B11").Value = "McGivney, Kluger, Clark & Intoccia, P.C."
$ws.Range("C11").Value = "Florham Park, NJ"
$ws.Range("D11").Value = "Full-Time"
$ws.Range("E11").Value = "We are seeking a Construction Project Engineer to join our team! Responsibilities: * Verify project drawing and specification compliance * Method of Procedure research and development * Construction ..."
$ws.Range("F11").Value = "https://www.ziprecruiter.com/k/l/AAIhcoIvzOu4Jh0BkUAZPGPXo-yJGs_6JmeR_sCioW5dn1lHl4sYV9qZhcUi6JVpyiXpiMzMzg-3XxDzLO8c1YcXBiygHPI2oR_A-OgbKgL9bNVT7y0mB8CAnYiaNOPcZVY4BHhNGupMIsPpXoRUlS9KE-zNdmRw1DY_sS2XFL_gVl0bw_Kov1M"
